$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column G (old G/H -> new I/J, etc.)
#    This shifts all existing data in columns G..U to I..W automatically,
#    carrying along per-column width/bestFit formatting.
$ws.Range("G1:H1").EntireColumn.Insert()

# 2. New header cells for the inserted columns.
$ws.Range("G1").Value = "Corresponding Organizer"
$ws.Range("H1").Value = "Email"

# 3. Row 2 (F1TENTH entry): add corresponding organizer + email, and swap
#    the now-duplicated "Organizer1" slot (I2) so Rahul Mangharam (the
#    corresponding organizer) leads the organizer list, with Hongrui Zheng
#    sliding into the slot Rahul previously occupied (O2).
$ws.Range("G2").Value = "Rahul Mangharam"
$ws.Range("H2").Value = "rahulm@seas.upenn.edu"
$ws.Range("I2").Value = "Rahul Mangharam"
$ws.Range("O2").Value = "Hongrui Zheng"

# 4. Row 4 (Robotic Grasping and Manipulation Competition / IROS On-Demand):
#    add livestream + date info and the long-form description.
$ws.Range("F4").Value = "October 20-23, 2020"
$ws.Range("E4").Value = "Livestream presentation November 6, 2020 9-10:30am EST"
$ws.Range("W4").Value = "The competition has two tracks: service robot track and manufacturing track. The service robot track has one simple task -- make five cups of iced Matcha green tea. There will be two manufacturing sub-tasks in the competition, disassembly and assembly using a NIST Task Board (NTB)."
$ws.Range("W4").Font.Color = 0

# 5. Widen the two newly inserted columns to fit their new header/content.
$ws.Columns.Item(15).ColumnWidth = 16.166666666666668
$ws.Columns.Item(16).ColumnWidth = 29.498697916666668

# 6. Hyperlink the corresponding organizer's email address.
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:rahulm@seas.upenn.edu", "", "", "mailto:rahulm@seas.upenn.edu")

# 7. Put the active selection where the editor last left it.
$ws.Range("W18").Select()
